# Apply the "spider_dropper_parts" workbook update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# --- Title change (A1) ---
$ws.Range("A1").Value = "Stupidly Simple Spider Dropper Parts"

# --- Row 6: "elastic string" -> "fishing wire for hanging", new price/qty/link ---
$ws.Range("A6").Value = "fishing wire for hanging"
$ws.Range("B6").Value = 6.29
$ws.Range("C6").Value = 656
$ws.Range("O6").Value = "https://www.amazon.com/gp/product/B08KZPHDPY/"

# --- Row 17: JGY 370 motor price/qty update, new Amazon link, old link moved to P17 ---
$ws.Range("P17").Value = "https://www.amazon.com/gp/product/B099JZ351N/"
$ws.Range("B17").Value = 18.99
$ws.Range("C17").Value = 2
$ws.Hyperlinks.Add($ws.Range("O17"), "https://www.amazon.com/dp/B0CZQMZDM8/")
$ws.Range("O17").Value = "https://www.amazon.com/dp/B0CZQMZDM8/"
$ws.Range("O17").Style = "Hyperlink"

# --- Row 23: PCB quantity needed 10 -> 5 ---
$ws.Range("C23").Value = 5

# --- New notes appended after row 47 ---
$ws.Range("A49").Value = "In bulk, cable ties are essentially free."
$ws.Range("A51").Value = "The DC version (including the Slightly Smarter upgrade) requires a 12VDC power supply."
$ws.Range("A53").Value = "The power connector and screw terminal could be omitted in favor of directly wiring pigtail power connectors to the PCB."

# --- Update selection to match the saved view state ---
$ws.Range("A54").Select()
